# feat: add 2022-Q4 data
#
# - The existing "2022-Q3" sheet (fund-level detail table) is cloned into a
#   brand-new sheet placed right after it; the clone keeps the name
#   "2022-Q3" and keeps the original fund data + formatting untouched.
# - The original "2022-Q3" sheet is renamed to "2022-Q4" and its contents
#   are replaced with the new quarter's fund detail table.
# - On the "总计" (totals) sheet, the 2022-Q3 summary row is turned into the
#   2022-Q4 summary row, and a new row is appended below with the original
#   2022-Q3 summary figures.

$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item(1)
$q3Sheet = $wb.Worksheets.Item(2)

$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# 1) Clone the current "2022-Q3" sheet (formatting included) into a new
#    sheet placed right after it - this clone becomes the permanent home
#    for the old Q3 fund data once sheet #2 itself becomes "2022-Q4".
# ---------------------------------------------------------------------
$q3Clone = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $q3Sheet)
$q3Sheet.UsedRange.Copy($q3Clone.Range("A1")) | Out-Null

$q3Sheet.Name = "2022-Q4"
$q3Clone.Name = "2022-Q3"
$q4Sheet = $q3Sheet

# ---------------------------------------------------------------------
# 2) Wipe sheet #2 (now "2022-Q4") and fill it in with the new quarter's
#    fund detail table.
# ---------------------------------------------------------------------
$q4Sheet.Cells.Clear() | Out-Null

$headers = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
for ($col = 2; $col -le 8; $col++) {
    $cell = $q4Sheet.Cells.Item(1, $col)
    $cell.Value = "'" + $headers[$col - 2]
    $cell.Style = "Normal"
}
# Match the header cell styling used elsewhere in the workbook (bold font +
# thin border cellXf), by copying the format from the totals sheet header.
$totalSheet.Range("B1").Copy() | Out-Null
$q4Sheet.Range("B1:H1").PasteSpecial($xlPasteFormats) | Out-Null

# code, name, scale, total-position, position-ratio, market-value, rank
$q4Rows = @(
    @("001404", "招商移动互联网产业股票A", "13.83", "90.58", "4.54", "0.6279", 10),
    @("015773", "招商移动互联网产业股票C", "8.53", "90.58", "4.54", "0.3873", 10),
    @("213008", "宝盈资源优选混合", "8.72", "89.15", "4.02", "0.3505", 10),
    @("008655", "招商科技创新混合A", "3.94", "93.96", "4.98", "0.1962", 4),
    @("008656", "招商科技创新混合C", "3.59", "93.96", "4.98", "0.1788", 4),
    @("006025", "诺安优化配置混合", "0.77", "91.43", "7.52", "0.0579", 8),
    @("501073", "华安智联混合（LOF）A", "3.27", "36.94", "1.32", "0.0432", 5),
    @("016380", "华宝专精特新混合A", "0.10", "90.16", "5.24", "0.0052", 2),
    @("016381", "华宝专精特新混合C", "0.01", "90.16", "5.24", "0.0005", 2),
    @("016071", "华安智联混合（LOF）C", "0.00", "36.94", "1.32", 0, 5)
)

$r = 2
foreach ($row in $q4Rows) {
    $q4Sheet.Cells.Item($r, 1).Value = $r - 2

    $q4Sheet.Cells.Item($r, 2).Value = "'" + $row[0]
    $q4Sheet.Cells.Item($r, 2).Style = "Normal"
    $q4Sheet.Cells.Item($r, 3).Value = "'" + $row[1]
    $q4Sheet.Cells.Item($r, 3).Style = "Normal"
    $q4Sheet.Cells.Item($r, 4).Value = "'" + $row[2]
    $q4Sheet.Cells.Item($r, 4).Style = "Normal"
    $q4Sheet.Cells.Item($r, 5).Value = "'" + $row[3]
    $q4Sheet.Cells.Item($r, 5).Style = "Normal"
    $q4Sheet.Cells.Item($r, 6).Value = "'" + $row[4]
    $q4Sheet.Cells.Item($r, 6).Style = "Normal"

    $marketValue = $row[5]
    if ($marketValue -eq 0) {
        $q4Sheet.Cells.Item($r, 7).Value = 0
    } else {
        $q4Sheet.Cells.Item($r, 7).Value = "'" + $marketValue
        $q4Sheet.Cells.Item($r, 7).Style = "Normal"
    }

    $q4Sheet.Cells.Item($r, 8).Value = $row[6]

    $r++
}

# Re-apply the column-A index-cell style (thin border / bold) used
# throughout the workbook.
$totalSheet.Range("A2").Copy() | Out-Null
$q4Sheet.Range("A2:A11").PasteSpecial($xlPasteFormats) | Out-Null

# ---------------------------------------------------------------------
# 3) Update the "总计" sheet: append the original 2022-Q3 summary as row 3,
#    then turn row 2 into the 2022-Q4 summary.
# ---------------------------------------------------------------------
$totalSheet.Cells.Item(3, 1).Value = 1
$totalSheet.Range("A2").Copy() | Out-Null
$totalSheet.Range("A3").PasteSpecial($xlPasteFormats) | Out-Null

$totalSheet.Cells.Item(3, 2).Value = "'2022-Q3"
$totalSheet.Cells.Item(3, 2).Style = "Normal"
$totalSheet.Cells.Item(3, 3).Value = 17
$totalSheet.Cells.Item(3, 4).Value = 2.32

$totalSheet.Cells.Item(2, 2).Value = "'2022-Q4"
$totalSheet.Cells.Item(2, 2).Style = "Normal"
$totalSheet.Cells.Item(2, 3).Value = 10
$totalSheet.Cells.Item(2, 4).Value = 1.85
